$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.117.91'
$ws.Range('E2').Value = '  +4.93%  '

$ws.Range('D3').Value = '3.785.00'
$ws.Range('E3').Value = '  +7.82%  '

$ws.Range('D4').Value = '''0.999'
$ws.Range('E4').Value = '  -0.68%  '

$ws.Range('D5').Value = '''427.06'
$ws.Range('E5').Value = '  +10.40%  '

$ws.Range('D6').Value = '''139.49'
$ws.Range('E6').Value = '  +15.58%  '

$ws.Range('D7').Value = '''0.613'
$ws.Range('E7').Value = '  +6.13%  '

$ws.Range('E8').Value = '  -0.23%  '

$ws.Range('D9').Value = '''0.741'
$ws.Range('E9').Value = '  +11.04%  '

$ws.Range('E10').Value = '  +6.87%  '

$ws.Range('D11').Value = '''0.0000328'
$ws.Range('E11').Value = '  +7.05%  '

$ws.Range('D12').Value = '''43.41'
$ws.Range('E12').Value = '  +14.14%  '

$ws.Range('D13').Value = '''10.70'
$ws.Range('E13').Value = '  +19.33%  '

$ws.Range('D14').Value = '4.382.50'
$ws.Range('E14').Value = '  +7.12%  '

$ws.Range('D15').Value = '''15.02'
$ws.Range('E15').Value = '  +19.14%  '

$ws.Range('E16').Value = '  +1.46%  '

$ws.Range('D17').Value = '3.793.76'
$ws.Range('E17').Value = '  +7.93%  '

$ws.Range('D18').Value = '''20.21'
$ws.Range('E18').Value = '  +9.69%  '

$ws.Range('E19').Value = '  +13.38%  '

$ws.Range('D20').Value = '66.221.40'
$ws.Range('E20').Value = '  +4.62%  '

$ws.Range('D21').Value = '''412.74'
$ws.Range('E21').Value = '  +6.90%  '

$ws.Range('D22').Value = '''15.28'
$ws.Range('E22').Value = '  +11.22%  '

$ws.Range('D23').Value = '''3.32'
$ws.Range('E23').Value = '  +17.33%  '

$ws.Range('D24').Value = '''85.78'
$ws.Range('E24').Value = '  +7.21%  '

$ws.Range('D25').Value = '''37.18'
$ws.Range('E25').Value = '  +12.47%  '

$ws.Range('D26').Value = '''9.87'
$ws.Range('E26').Value = '  +48.95%  '

$ws.Range('E27').Value = '  +12.04%  '

$ws.Range('D28').Value = '''9.86'
$ws.Range('E28').Value = '  +15.81%  '

$ws.Range('E29').Value = '  -1.10%  '

$ws.Range('D30').Value = '''13.98'
$ws.Range('E30').Value = '  +20.71%  '

$ws.Range('D31').Value = '''709.19'
$ws.Range('E31').Value = '  +6.99%  '

$ws.Range('E32').Value = '  +19.09%  '

$ws.Range('D33').Value = '''2.78'
$ws.Range('E33').Value = '  +4.82%  '

$ws.Range('D34').Value = '''40.34'
$ws.Range('E34').Value = '  +12.17%  '

$ws.Range('D35').Value = '''0.999'
$ws.Range('E35').Value = '  -0.08%  '

$ws.Range('D36').Value = '''5.85'
$ws.Range('E36').Value = '  +45.06%  '

$ws.Range('D37').Value = '''0.151'
$ws.Range('E37').Value = '  +2.66%  '

$ws.Range('D38').Value = '''56.19'
$ws.Range('E38').Value = '  +5.90%  '

$ws.Range('D39').Value = '''0.0475'
$ws.Range('E39').Value = '  +10.56%  '

$ws.Range('B40').Value = 'PEPE'
$ws.Range('C40').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D40').Value = '0.0₃0684'
$ws.Range('E40').Value = '  +14.91%  '

$ws.Range('B41').Value = 'Fetch.AI'
$ws.Range('C41').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D41').Value = '''2.61'
$ws.Range('E41').Value = '  +51.96%  '

$ws.Range('E42').Value = '  +9.54%  '

$ws.Range('E43').Value = '  +0.08%  '

$ws.Range('E44').Value = '  +9.03%  '

$ws.Range('D45').Value = '''3.39'
$ws.Range('E45').Value = '  +12.50%  '

$ws.Range('D46').Value = '''0.322'
$ws.Range('E46').Value = '  +19.75%  '

$ws.Range('D47').Value = '''3.16'
$ws.Range('E47').Value = '  +4.62%  '

$ws.Range('E48').Value = '  +8.08%  '

$ws.Range('D49').Value = '''2.65'
$ws.Range('E49').Value = '  +7.13%  '

$ws.Range('D50').Value = '''142.61'
$ws.Range('E50').Value = '  +2.94%  '

$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''25.91'
$ws.Range('E51').Value = '  +0.61%  '
